$d = $word.ActiveDocument

# Common run properties shared by all new paragraphs / runs being added.
$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="1D1C1D"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="tr-TR"/></w:rPr>'

# Paragraph 1: a completely empty paragraph (no runs), matching the blank
# line that separates the previous answer from the next question.
$p1 = '<w:p><w:pPr><w:jc w:val="both"/>' + $rPr + '</w:pPr></w:p>'

# Paragraph 2: "Q: What are the metrics used to evaluate a Regression Model?"
$p2 = '<w:p><w:pPr><w:jc w:val="both"/>' + $rPr + '</w:pPr>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve">Q: </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>What</w:t></w:r><w:proofErr w:type="spellEnd"/>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>are</w:t></w:r><w:proofErr w:type="spellEnd"/>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>the</w:t></w:r><w:proofErr w:type="spellEnd"/>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>metrics</w:t></w:r><w:proofErr w:type="spellEnd"/>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>used</w:t></w:r><w:proofErr w:type="spellEnd"/>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>to</w:t></w:r><w:proofErr w:type="spellEnd"/>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>evaluate</w:t></w:r><w:proofErr w:type="spellEnd"/>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve"> a </w:t></w:r>' `
    + '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>Regression</w:t></w:r><w:proofErr w:type="spellEnd"/>' `
    + '<w:r>' + $rPr + '<w:t xml:space="preserve"> Model?</w:t></w:r>' `
    + '</w:p>'

# Paragraph 3: "A1:"
$p3 = '<w:p><w:pPr><w:jc w:val="both"/>' + $rPr + '</w:pPr>' `
    + '<w:r>' + $rPr + '<w:t>A1:</w:t></w:r>' `
    + '</w:p>'

$bodyFragment = $p1 + $p2 + $p3

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' `
    + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
    + '<pkg:xmlData>' `
    + '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
    + '<w:body>' + $bodyFragment + '</w:body></w:document>' `
    + '</pkg:xmlData></pkg:part></pkg:package>'

# Collapse to the very end of the document (after the last paragraph,
# "...underfitting the data.") and insert the new paragraphs there, right
# before the closing sectPr.
$r = $d.Content
$r.Collapse(0)
$r.InsertXML($xml)
